$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"
$errorText = "The handback type mt is not match with handoff type ht."

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Cells.Item(6,3).Value = $statusText
$wsZh.Cells.Item(6,11).Value = $errorText

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Cells.Item(6,3).Value = $statusText
$wsDe.Cells.Item(6,11).Value = $errorText
